# Append a new book ("The Leadership Challenge") to the end of the
# "Completed" reading list on sheet1 (row 38).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 38

$ws.Cells.Item($row, 1).Value = "The Leadership Challenge"
$ws.Cells.Item($row, 2).Value = "Barry Posner;James Kouzes"

# Copy the number format (m/d/yyyy style) used by the existing date
# columns onto the new row before filling in the values, so the new
# cells reuse the same cell style as the rest of the column instead of
# Excel inventing a brand-new style.
$ws.Range("C37:D37").Copy()
$ws.Range("C" + $row + ":D" + $row).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 3).Value = 43904
$ws.Cells.Item($row, 4).Value = 43906
$ws.Cells.Item($row, 5).Value = "leadership;success;self-improvement"
$ws.Cells.Item($row, 6).Value = "Audio"
$ws.Cells.Item($row, 7).Value = "11 Hours 1 Min"

# Mirror the view-state change that happens in Excel when the new row
# is entered: the active cell moves on to the next empty row.
$ws.Range("A39").Select()
$excel.ActiveWindow.ScrollRow = 19

